$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    4  = 127836377946
    5  = 10260.2478
    6  = 344.180383
    7  = 117.858539
    8  = 66.853517
    9  = 510.82979
    10 = 121.271979
    11 = 1598.4183
    12 = 927.39034
    13 = 491.459075
    14 = 82172205.13699999
    15 = 8596.4746
    16 = 152.204531
    17 = 4323.02753
    18 = 474.12874
    19 = 130.203457
    20 = 3713.7039
    21 = 4044.47962
    22 = 244.75483
    23 = 251.45151
    24 = 1175.86616
    25 = 491.98218
    26 = 388.93089
    27 = 867.29317
    28 = 1338.37435
    29 = 504.64666
    30 = 661.22308
    31 = 491.39772
    32 = 4232030.437000001
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
